# Add two new columns, I ("I0") and J ("IF"), to Sheet1, with per-row
# computed values, matching the committed change "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - same bordered/bold/centered style as the other
# header cells (e.g. H1), applied by copying the existing header style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row I and J values (rows 2-58).
$iVals = @(6,7,7,6,7,10,7,7,9,7,6,7,7,7,8,7,5,7,7,8,8,5,6,5,6,7,6,5,6,7,7,6,8,8,7,7,7,4,7,7,8,8,8,7,7,9,6,9,5,8,5,5,4,8,8,3,4)
$jVals = @(6,7,7,7,7,10,7,7,9,8,6,7,8,8,9,7,5,7,7,8,8,5,7,5,7,8,6,6,6,8,7,6,8,8,7,8,7,4,7,8,8,9,8,8,8,9,7,9,6,9,5,5,5,8,8,3,4)

for ($r = 2; $r -le 58; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]   # column I
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]  # column J
}
